# Add a new "UK" market test-data sheet, cloned from the existing "Poland"
# sheet (same layout/styling), with UK-specific values and two extra rows.

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Clone the Poland sheet (keeps column widths, styles, merged cells, etc.)
# and place the copy immediately after Poland - this also makes it the
# active/selected sheet, matching how Excel behaves when you duplicate a
# tab via the UI.
$poland.Copy([System.Reflection.Missing]::Value, $poland) | Out-Null

$uk = $wb.Worksheets.Item("Poland (2)")
$uk.Name = "UK"

# Country-specific ticket / market values (set B4 before B2 so the shared
# strings land in the same order as the source edit: ticket code first,
# then the market label).
$uk.Range("B4").Value = "NGC-2741/T3366"
$uk.Range("B2").Value = "UK Market"

# Extend the repeaters list by two rows, reusing row 17's border styling.
$uk.Range("A17").Copy() | Out-Null
$uk.Range("A18:A19").PasteSpecial(-4122) | Out-Null
$uk.Range("A18").Value = "Wg"
$uk.Range("A19").Value = "Repeaters"

# Match the author's final selection on the new sheet.
$uk.Range("B4").Select() | Out-Null
